$d = $word.ActiveDocument

$d.Content.Find.Execute('2024-11-18 Monday', $true, $false, $false, $false, $false, $true, 1, $false, '2024-11-19 Tuesday', 2) | Out-Null
$d.Content.Find.Execute('414÷2=', $true, $false, $false, $false, $false, $true, 1, $false, '719÷4=', 2) | Out-Null
$d.Content.Find.Execute('488÷4=', $true, $false, $false, $false, $false, $true, 1, $false, '920÷2=', 2) | Out-Null
$d.Content.Find.Execute('519÷3=', $true, $false, $false, $false, $false, $true, 1, $false, '112÷9=', 2) | Out-Null
$d.Content.Find.Execute('881÷3=', $true, $false, $false, $false, $false, $true, 1, $false, '562÷5=', 2) | Out-Null
$d.Content.Find.Execute('250÷6=', $true, $false, $false, $false, $false, $true, 1, $false, '929÷3=', 2) | Out-Null
$d.Content.Find.Execute('612÷2=', $true, $false, $false, $false, $false, $true, 1, $false, '468÷7=', 2) | Out-Null
$d.Content.Find.Execute('868÷9=', $true, $false, $false, $false, $false, $true, 1, $false, '746÷5=', 2) | Out-Null
$d.Content.Find.Execute('810÷7=', $true, $false, $false, $false, $false, $true, 1, $false, '279÷5=', 2) | Out-Null
$d.Content.Find.Execute('134÷3=', $true, $false, $false, $false, $false, $true, 1, $false, '522÷5=', 2) | Out-Null
$d.Content.Find.Execute('526÷6=', $true, $false, $false, $false, $false, $true, 1, $false, '483÷5=', 2) | Out-Null
$d.Content.Find.Execute('665÷3=', $true, $false, $false, $false, $false, $true, 1, $false, '355÷5=', 2) | Out-Null
$d.Content.Find.Execute('642÷3=', $true, $false, $false, $false, $false, $true, 1, $false, '546÷8=', 2) | Out-Null
$d.Content.Find.Execute('708÷3=', $true, $false, $false, $false, $false, $true, 1, $false, '654÷6=', 2) | Out-Null
$d.Content.Find.Execute('498÷6=', $true, $false, $false, $false, $false, $true, 1, $false, '462÷9=', 2) | Out-Null
$d.Content.Find.Execute('822÷4=', $true, $false, $false, $false, $false, $true, 1, $false, '935÷9=', 2) | Out-Null
$d.Content.Find.Execute('955÷2=', $true, $false, $false, $false, $false, $true, 1, $false, '303÷7=', 2) | Out-Null
$d.Content.Find.Execute('864÷8=', $true, $false, $false, $false, $false, $true, 1, $false, '342÷9=', 2) | Out-Null
$d.Content.Find.Execute('129÷2=', $true, $false, $false, $false, $false, $true, 1, $false, '135÷5=', 2) | Out-Null
$d.Content.Find.Execute('656÷5=', $true, $false, $false, $false, $false, $true, 1, $false, '586÷8=', 2) | Out-Null
$d.Content.Find.Execute('736÷4=', $true, $false, $false, $false, $false, $true, 1, $false, '659÷5=', 2) | Out-Null
$d.Content.Find.Execute('782÷6=', $true, $false, $false, $false, $false, $true, 1, $false, '205÷9=', 2) | Out-Null
$d.Content.Find.Execute('132÷5=', $true, $false, $false, $false, $false, $true, 1, $false, '942÷8=', 2) | Out-Null
$d.Content.Find.Execute('334÷9=', $true, $false, $false, $false, $false, $true, 1, $false, '569÷2=', 2) | Out-Null
$d.Content.Find.Execute('337÷2=', $true, $false, $false, $false, $false, $true, 1, $false, '830÷8=', 2) | Out-Null
$d.Content.Find.Execute('889÷2=', $true, $false, $false, $false, $false, $true, 1, $false, '804÷4=', 2) | Out-Null
